$wb = $excel.ActiveWorkbook

# Add the new "1017 nodes" sheet before the current first sheet
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "1017 nodes"

# Populate cells (A2:H44) with the node data
$ws.Cells.Item(2,1).Value = 'Brain Region Name'
$ws.Cells.Item(2,2).Value = 'Right Hemisphere'
$ws.Cells.Item(2,4).Value = 'Left Hemisphere'
$ws.Cells.Item(2,6).Value = 'Lobe Code'
$ws.Cells.Item(2,7).Value = 'Lobe Name'
$ws.Cells.Item(3,1).Value = 'lateral orbitofrontal'
$ws.Cells.Item(3,2).Value = 1
$ws.Cells.Item(3,3).Value = 17
$ws.Cells.Item(3,4).Value = 509
$ws.Cells.Item(3,5).Value = 525
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 'Frontal'
$ws.Cells.Item(3,8).Value = 'Frontal'
$ws.Cells.Item(4,1).Value = 'parsobitalis'
$ws.Cells.Item(4,2).Value = 18
$ws.Cells.Item(4,3).Value = 21
$ws.Cells.Item(4,4).Value = 526
$ws.Cells.Item(4,5).Value = 529
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 'Frontal'
$ws.Cells.Item(4,8).Value = 'Frontal'
$ws.Cells.Item(5,1).Value = 'frontal pole'
$ws.Cells.Item(5,2).Value = 22
$ws.Cells.Item(5,3).Value = 23
$ws.Cells.Item(5,4).Value = 530
$ws.Cells.Item(5,5).Value = 531
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 'Frontal'
$ws.Cells.Item(5,8).Value = 'Frontal'
$ws.Cells.Item(6,1).Value = 'medial orbitofrontal'
$ws.Cells.Item(6,2).Value = 24
$ws.Cells.Item(6,3).Value = 34
$ws.Cells.Item(6,4).Value = 532
$ws.Cells.Item(6,5).Value = 542
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 'Frontal'
$ws.Cells.Item(6,8).Value = 'Frontal'
$ws.Cells.Item(7,1).Value = 'parstriangularis'
$ws.Cells.Item(7,2).Value = 35
$ws.Cells.Item(7,3).Value = 42
$ws.Cells.Item(7,4).Value = 543
$ws.Cells.Item(7,5).Value = 550
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 'Frontal'
$ws.Cells.Item(7,8).Value = 'Frontal'
$ws.Cells.Item(8,1).Value = 'parsopecularis'
$ws.Cells.Item(8,2).Value = 43
$ws.Cells.Item(8,3).Value = 51
$ws.Cells.Item(8,4).Value = 551
$ws.Cells.Item(8,5).Value = 559
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 'Frontal'
$ws.Cells.Item(8,8).Value = 'Frontal'
$ws.Cells.Item(9,1).Value = 'rostral middle frontal'
$ws.Cells.Item(9,2).Value = 52
$ws.Cells.Item(9,3).Value = 78
$ws.Cells.Item(9,4).Value = 560
$ws.Cells.Item(9,5).Value = 586
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 'Frontal'
$ws.Cells.Item(9,8).Value = 'Frontal'
$ws.Cells.Item(10,1).Value = 'superior frontal'
$ws.Cells.Item(10,2).Value = 79
$ws.Cells.Item(10,3).Value = 120
$ws.Cells.Item(10,4).Value = 587
$ws.Cells.Item(10,5).Value = 628
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 'Frontal'
$ws.Cells.Item(10,8).Value = 'Frontal'
$ws.Cells.Item(11,1).Value = 'caudal middle frontal'
$ws.Cells.Item(11,2).Value = 121
$ws.Cells.Item(11,3).Value = 131
$ws.Cells.Item(11,4).Value = 629
$ws.Cells.Item(11,5).Value = 639
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 'Frontal'
$ws.Cells.Item(11,8).Value = 'Frontal'
$ws.Cells.Item(12,1).Value = 'precentral'
$ws.Cells.Item(12,2).Value = 132
$ws.Cells.Item(12,3).Value = 167
$ws.Cells.Item(12,4).Value = 640
$ws.Cells.Item(12,5).Value = 675
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 'Frontal'
$ws.Cells.Item(12,8).Value = 'Frontal'
$ws.Cells.Item(13,1).Value = 'paracentral'
$ws.Cells.Item(13,2).Value = 168
$ws.Cells.Item(13,3).Value = 179
$ws.Cells.Item(13,4).Value = 676
$ws.Cells.Item(13,5).Value = 687
$ws.Cells.Item(13,6).Value = 2
$ws.Cells.Item(13,7).Value = 'Parietal'
$ws.Cells.Item(13,8).Value = 'Parietal'
$ws.Cells.Item(14,1).Value = 'rostral anterior cingulate'
$ws.Cells.Item(14,2).Value = 180
$ws.Cells.Item(14,3).Value = 183
$ws.Cells.Item(14,4).Value = 688
$ws.Cells.Item(14,5).Value = 691
$ws.Cells.Item(14,6).Value = 2
$ws.Cells.Item(14,7).Value = 'Parietal'
$ws.Cells.Item(14,8).Value = 'Limbic'
$ws.Cells.Item(15,1).Value = 'caudal anterior cingulate'
$ws.Cells.Item(15,2).Value = 184
$ws.Cells.Item(15,3).Value = 189
$ws.Cells.Item(15,4).Value = 692
$ws.Cells.Item(15,5).Value = 697
$ws.Cells.Item(15,6).Value = 2
$ws.Cells.Item(15,7).Value = 'Parietal'
$ws.Cells.Item(15,8).Value = 'Limbic'
$ws.Cells.Item(16,1).Value = 'posterior cingulate'
$ws.Cells.Item(16,2).Value = 190
$ws.Cells.Item(16,3).Value = 198
$ws.Cells.Item(16,4).Value = 698
$ws.Cells.Item(16,5).Value = 706
$ws.Cells.Item(16,6).Value = 2
$ws.Cells.Item(16,7).Value = 'Parietal'
$ws.Cells.Item(16,8).Value = 'Limbic'
$ws.Cells.Item(17,1).Value = 'Isthmus cingulate'
$ws.Cells.Item(17,2).Value = 199
$ws.Cells.Item(17,3).Value = 204
$ws.Cells.Item(17,4).Value = 707
$ws.Cells.Item(17,5).Value = 712
$ws.Cells.Item(17,6).Value = 2
$ws.Cells.Item(17,7).Value = 'Parietal'
$ws.Cells.Item(17,8).Value = 'Limbic'
$ws.Cells.Item(18,1).Value = 'postcentral'
$ws.Cells.Item(18,2).Value = 205
$ws.Cells.Item(18,3).Value = 234
$ws.Cells.Item(18,4).Value = 713
$ws.Cells.Item(18,5).Value = 742
$ws.Cells.Item(18,6).Value = 2
$ws.Cells.Item(18,7).Value = 'Parietal'
$ws.Cells.Item(18,8).Value = 'Parietal'
$ws.Cells.Item(19,1).Value = 'supramarginal'
$ws.Cells.Item(19,2).Value = 235
$ws.Cells.Item(19,3).Value = 254
$ws.Cells.Item(19,4).Value = 743
$ws.Cells.Item(19,5).Value = 762
$ws.Cells.Item(19,6).Value = 2
$ws.Cells.Item(19,7).Value = 'Parietal'
$ws.Cells.Item(19,8).Value = 'Parietal'
$ws.Cells.Item(20,1).Value = 'superior parietal'
$ws.Cells.Item(20,2).Value = 255
$ws.Cells.Item(20,3).Value = 283
$ws.Cells.Item(20,4).Value = 763
$ws.Cells.Item(20,5).Value = 791
$ws.Cells.Item(20,6).Value = 2
$ws.Cells.Item(20,7).Value = 'Parietal'
$ws.Cells.Item(20,8).Value = 'Parietal'
$ws.Cells.Item(21,1).Value = 'inferior parietal'
$ws.Cells.Item(21,2).Value = 284
$ws.Cells.Item(21,3).Value = 309
$ws.Cells.Item(21,4).Value = 792
$ws.Cells.Item(21,5).Value = 817
$ws.Cells.Item(21,6).Value = 2
$ws.Cells.Item(21,7).Value = 'Parietal'
$ws.Cells.Item(21,8).Value = 'Parietal'
$ws.Cells.Item(22,1).Value = 'precuneus'
$ws.Cells.Item(22,2).Value = 310
$ws.Cells.Item(22,3).Value = 332
$ws.Cells.Item(22,4).Value = 818
$ws.Cells.Item(22,5).Value = 840
$ws.Cells.Item(22,6).Value = 2
$ws.Cells.Item(22,7).Value = 'Parietal'
$ws.Cells.Item(22,8).Value = 'Parietal'
$ws.Cells.Item(23,1).Value = 'cuneus'
$ws.Cells.Item(23,2).Value = 333
$ws.Cells.Item(23,3).Value = 340
$ws.Cells.Item(23,4).Value = 841
$ws.Cells.Item(23,5).Value = 848
$ws.Cells.Item(23,6).Value = 3
$ws.Cells.Item(23,7).Value = 'Occiptal'
$ws.Cells.Item(23,8).Value = 'Occipital'
$ws.Cells.Item(24,1).Value = 'pericalcarine'
$ws.Cells.Item(24,2).Value = 341
$ws.Cells.Item(24,3).Value = 348
$ws.Cells.Item(24,4).Value = 849
$ws.Cells.Item(24,5).Value = 856
$ws.Cells.Item(24,6).Value = 3
$ws.Cells.Item(24,7).Value = 'Occiptal'
$ws.Cells.Item(24,8).Value = 'Occipital'
$ws.Cells.Item(25,1).Value = 'lateral occipital'
$ws.Cells.Item(25,2).Value = 349
$ws.Cells.Item(25,3).Value = 371
$ws.Cells.Item(25,4).Value = 857
$ws.Cells.Item(25,5).Value = 879
$ws.Cells.Item(25,6).Value = 3
$ws.Cells.Item(25,7).Value = 'Occiptal'
$ws.Cells.Item(25,8).Value = 'Occipital'
$ws.Cells.Item(26,1).Value = 'lingual'
$ws.Cells.Item(26,2).Value = 372
$ws.Cells.Item(26,3).Value = 388
$ws.Cells.Item(26,4).Value = 880
$ws.Cells.Item(26,5).Value = 896
$ws.Cells.Item(26,6).Value = 3
$ws.Cells.Item(26,7).Value = 'Occiptal'
$ws.Cells.Item(26,8).Value = 'Occipital'
$ws.Cells.Item(27,1).Value = 'fusiform'
$ws.Cells.Item(27,2).Value = 389
$ws.Cells.Item(27,3).Value = 405
$ws.Cells.Item(27,4).Value = 897
$ws.Cells.Item(27,5).Value = 913
$ws.Cells.Item(27,6).Value = 4
$ws.Cells.Item(27,7).Value = 'Temporal'
$ws.Cells.Item(27,8).Value = 'Temporal'
$ws.Cells.Item(28,1).Value = 'parahippocampal'
$ws.Cells.Item(28,2).Value = 406
$ws.Cells.Item(28,3).Value = 411
$ws.Cells.Item(28,4).Value = 914
$ws.Cells.Item(28,5).Value = 919
$ws.Cells.Item(28,6).Value = 4
$ws.Cells.Item(28,7).Value = 'Temporal'
$ws.Cells.Item(28,8).Value = 'Limbic'
$ws.Cells.Item(29,1).Value = 'entorhinal'
$ws.Cells.Item(29,2).Value = 412
$ws.Cells.Item(29,3).Value = 413
$ws.Cells.Item(29,4).Value = 920
$ws.Cells.Item(29,5).Value = 921
$ws.Cells.Item(29,6).Value = 4
$ws.Cells.Item(29,7).Value = 'Temporal'
$ws.Cells.Item(29,8).Value = 'Limbic'
$ws.Cells.Item(30,1).Value = 'temporal pole'
$ws.Cells.Item(30,2).Value = 414
$ws.Cells.Item(30,3).Value = 416
$ws.Cells.Item(30,4).Value = 922
$ws.Cells.Item(30,5).Value = 924
$ws.Cells.Item(30,6).Value = 4
$ws.Cells.Item(30,7).Value = 'Temporal'
$ws.Cells.Item(30,8).Value = 'Temporal'
$ws.Cells.Item(31,1).Value = 'inferior temporal'
$ws.Cells.Item(31,2).Value = 417
$ws.Cells.Item(31,3).Value = 432
$ws.Cells.Item(31,4).Value = 925
$ws.Cells.Item(31,5).Value = 940
$ws.Cells.Item(31,6).Value = 4
$ws.Cells.Item(31,7).Value = 'Temporal'
$ws.Cells.Item(31,8).Value = 'Temporal'
$ws.Cells.Item(32,1).Value = 'middle temporal'
$ws.Cells.Item(32,2).Value = 433
$ws.Cells.Item(32,3).Value = 451
$ws.Cells.Item(32,4).Value = 941
$ws.Cells.Item(32,5).Value = 959
$ws.Cells.Item(32,6).Value = 4
$ws.Cells.Item(32,7).Value = 'Temporal'
$ws.Cells.Item(32,8).Value = 'Temporal'
$ws.Cells.Item(33,1).Value = 'bankssts'
$ws.Cells.Item(33,2).Value = 452
$ws.Cells.Item(33,3).Value = 457
$ws.Cells.Item(33,4).Value = 960
$ws.Cells.Item(33,5).Value = 965
$ws.Cells.Item(33,6).Value = 4
$ws.Cells.Item(33,7).Value = 'Temporal'
$ws.Cells.Item(33,8).Value = 'Temporal'
$ws.Cells.Item(34,1).Value = 'superior temporal'
$ws.Cells.Item(34,2).Value = 458
$ws.Cells.Item(34,3).Value = 482
$ws.Cells.Item(34,4).Value = 966
$ws.Cells.Item(34,5).Value = 990
$ws.Cells.Item(34,6).Value = 4
$ws.Cells.Item(34,7).Value = 'Temporal'
$ws.Cells.Item(34,8).Value = 'Temporal'
$ws.Cells.Item(35,1).Value = 'transverse temporal'
$ws.Cells.Item(35,2).Value = 483
$ws.Cells.Item(35,3).Value = 485
$ws.Cells.Item(35,4).Value = 991
$ws.Cells.Item(35,5).Value = 993
$ws.Cells.Item(35,6).Value = 4
$ws.Cells.Item(35,7).Value = 'Temporal'
$ws.Cells.Item(35,8).Value = 'Temporal'
$ws.Cells.Item(36,1).Value = 'insula'
$ws.Cells.Item(36,2).Value = 486
$ws.Cells.Item(36,3).Value = 501
$ws.Cells.Item(36,4).Value = 994
$ws.Cells.Item(36,5).Value = 1009
$ws.Cells.Item(36,6).Value = 4
$ws.Cells.Item(36,7).Value = 'Temporal'
$ws.Cells.Item(36,8).Value = 'Temporal'
$ws.Cells.Item(37,1).Value = 'thalamus proper'
$ws.Cells.Item(37,2).Value = 502
$ws.Cells.Item(37,3).Value = 502
$ws.Cells.Item(37,4).Value = 1010
$ws.Cells.Item(37,5).Value = 1010
$ws.Cells.Item(37,6).Value = 5
$ws.Cells.Item(37,7).Value = 'Subcortical'
$ws.Cells.Item(37,8).Value = 'Basal Ganglia'
$ws.Cells.Item(38,1).Value = 'caudate'
$ws.Cells.Item(38,2).Value = 503
$ws.Cells.Item(38,3).Value = 503
$ws.Cells.Item(38,4).Value = 1011
$ws.Cells.Item(38,5).Value = 1011
$ws.Cells.Item(38,6).Value = 5
$ws.Cells.Item(38,7).Value = 'Subcortical'
$ws.Cells.Item(38,8).Value = 'Basal Ganglia'
$ws.Cells.Item(39,1).Value = 'putamen'
$ws.Cells.Item(39,2).Value = 504
$ws.Cells.Item(39,3).Value = 504
$ws.Cells.Item(39,4).Value = 1012
$ws.Cells.Item(39,5).Value = 1012
$ws.Cells.Item(39,6).Value = 5
$ws.Cells.Item(39,7).Value = 'Subcortical'
$ws.Cells.Item(39,8).Value = 'Basal Ganglia'
$ws.Cells.Item(40,1).Value = 'pallidum'
$ws.Cells.Item(40,2).Value = 505
$ws.Cells.Item(40,3).Value = 505
$ws.Cells.Item(40,4).Value = 1013
$ws.Cells.Item(40,5).Value = 1013
$ws.Cells.Item(40,6).Value = 5
$ws.Cells.Item(40,7).Value = 'Subcortical'
$ws.Cells.Item(40,8).Value = 'Basal Ganglia'
$ws.Cells.Item(41,1).Value = 'accumbens area'
$ws.Cells.Item(41,2).Value = 506
$ws.Cells.Item(41,3).Value = 506
$ws.Cells.Item(41,4).Value = 1014
$ws.Cells.Item(41,5).Value = 1014
$ws.Cells.Item(41,6).Value = 5
$ws.Cells.Item(41,7).Value = 'Subcortical'
$ws.Cells.Item(41,8).Value = 'Basal Ganglia'
$ws.Cells.Item(42,1).Value = 'hippocampus'
$ws.Cells.Item(42,2).Value = 507
$ws.Cells.Item(42,3).Value = 507
$ws.Cells.Item(42,4).Value = 1015
$ws.Cells.Item(42,5).Value = 1015
$ws.Cells.Item(42,6).Value = 5
$ws.Cells.Item(42,7).Value = 'Subcortical'
$ws.Cells.Item(42,8).Value = 'Temporal'
$ws.Cells.Item(43,1).Value = 'amygdala'
$ws.Cells.Item(43,2).Value = 508
$ws.Cells.Item(43,3).Value = 508
$ws.Cells.Item(43,4).Value = 1016
$ws.Cells.Item(43,5).Value = 1016
$ws.Cells.Item(43,6).Value = 5
$ws.Cells.Item(43,7).Value = 'Subcortical'
$ws.Cells.Item(43,8).Value = 'Basal Ganglia'
$ws.Cells.Item(44,1).Value = 'brainstem'
$ws.Cells.Item(44,4).Value = 1017
$ws.Cells.Item(44,5).Value = 1017
$ws.Cells.Item(44,6).Value = 5
$ws.Cells.Item(44,7).Value = 'Subcortical'
$ws.Cells.Item(44,8).Value = 'Brainstem'

# Match the saved view state: active cell E14 selected, this sheet active/tabSelected
$ws.Range("E14").Select()
